$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.919.95'
$ws.Range('E2').Value = '  -0.93%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.13'
$ws.Range('E3').Value = '  -0.50%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.92'
$ws.Range('E5').Value = '  -0.78%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5047'
$ws.Range('E7').Value = '  -1.43%  '

# Row 8
$ws.Range('E8').Value = '  -2.75%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07169'
$ws.Range('E9').Value = '  +0.41%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8944'
$ws.Range('E10').Value = '  +0.81%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.77'
$ws.Range('E11').Value = '  +0.69%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.871.16'
$ws.Range('E12').Value = '  -0.04%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07478'
$ws.Range('E13').Value = '  -0.89%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.62'
$ws.Range('E14').Value = '  +3.80%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.228'
$ws.Range('E15').Value = '  -1.83%  '

# Row 16
$ws.Range('E16').Value = '  -0.03%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008496'
$ws.Range('E17').Value = '  +0.29%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.19'
$ws.Range('E18').Value = '  +0.81%  '

# Row 19
$ws.Range('E19').Value = '  -0.06%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.951.59'
$ws.Range('E20').Value = '  -0.97%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.035'
$ws.Range('E21').Value = '  -0.43%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.080.17'
$ws.Range('E22').Value = '  -1.71%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.39'
$ws.Range('E23').Value = '  -1.65%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.397'
$ws.Range('E24').Value = '  -1.21%  '

# Row 25
$ws.Range('E25').Value = '  -1.91%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.790'
$ws.Range('E26').Value = '  -3.10%  '

# Row 27
$ws.Range('E27').Value = '  -0.44%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.081'
$ws.Range('E28').Value = '  -0.88%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.13'
$ws.Range('E29').Value = '  +0.13%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.703'
$ws.Range('E30').Value = '  -0.34%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.676'
$ws.Range('E31').Value = '  +0.10%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09247'
$ws.Range('E32').Value = '  +2.41%  '

# Row 33
$ws.Range('E33').Value = '  -0.79%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7529'
$ws.Range('E34').Value = '  +2.38%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.993'
$ws.Range('E35').Value = '  -3.02%  '

# Row 36
$ws.Range('E36').Value = '  -0.52%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.263'
$ws.Range('E37').Value = '  +6.64%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.535'
$ws.Range('E38').Value = '  +1.09%  '

# Row 39
$ws.Range('E39').Value = '  -2.45%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5588'
$ws.Range('E40').Value = '  +4.15%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.070'
$ws.Range('E41').Value = '  -0.34%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '118.90'
$ws.Range('E42').Value = '  +1.55%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.533'
$ws.Range('E43').Value = '  -0.67%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.529'
$ws.Range('E44').Value = '  +2.40%  '

# Row 45
$ws.Range('E45').Value = '  -0.16%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4697'
$ws.Range('E46').Value = '  +1.21%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9997'
$ws.Range('E47').Value = '  -0.02%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.08'
$ws.Range('E48').Value = '  +0.24%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.563'
$ws.Range('E49').Value = '  -0.47%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.80'
$ws.Range('E50').Value = '  +0.80%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.91'
$ws.Range('E51').Value = '  -2.28%  '
